$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format (style) from the last populated row (A17) down onto
# the new date cells A18:A20 before writing values into them, so they pick
# up the same numFmtId-22 date/time style instead of Excel inventing a new
# custom numFmt.
$ws.Range("A17").Copy()
$ws.Range("A18:A20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 18 - "File reorg, gamemode system"
$ws.Range("A18").Value = 45256.945138888892
$ws.Range("B18").Value = 1786880
$ws.Range("C18").Value = 693248
$ws.Range("D18").Value = 396288
$ws.Range("F18").Value = 168915701
$ws.Range("I18").Value = "File reorg, gamemode system"

# Row 19 - "Remove cinematics"
$ws.Range("A19").Value = 45257.870833333334
$ws.Range("B19").Value = 1776640
$ws.Range("C19").Value = 693248
$ws.Range("D19").Value = 396288
$ws.Range("F19").Value = 168970214
$ws.Range("I19").Value = "Remove cinematics"

# Row 20 - "Remove now-redundant menu options, remaining PCX support"
$ws.Range("A20").Value = 45257.890277777777
$ws.Range("B20").Value = 1776128
$ws.Range("C20").Value = 693248
$ws.Range("D20").Value = 392704
$ws.Range("F20").Value = 168986319
$ws.Range("I20").Value = "Remove now-redundant menu options, remaining PCX support"

# Move the active selection to I29, matching the author's cursor position
# after entering the new rows.
$ws.Range("I29").Select() | Out-Null
